# Generate Report for Handback
#
# The nightly CI run has handed the localized files back ("in sync with
# en-US") for both language sheets (zh-cn, de-de). This updates the
# localization-status report:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The (previously empty) "Latest Target File" / "Latest Handback File"
#     columns get populated with hyperlinked file names
#   - "Latest Handback DateTime" gets a real timestamp instead of the
#     zero-date placeholder

$wb = $excel.ActiveWorkbook

# Cornflower blue (0x6495ED) packed the way Excel's Font.Color (a BGR Long) expects.
$hyperlinkColor = 15570276

function Set-HandbackRow($SheetName, $Row, $HandbackDateTime, $SourceMdName, $SourceMdUrl, $XlfName, $XlfUrl) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Cells.Item($Row, 3).Value = "Handed back: in sync with en-US"

    # Latest Target File (F) - the file handed back, same name as the source file
    $fCell = $ws.Cells.Item($Row, 6)
    $fCell.Value = $SourceMdName
    $ws.Hyperlinks.Add($fCell, $SourceMdUrl, "", "", $SourceMdName) | Out-Null
    $fCell.Font.Underline = 2
    $fCell.Font.Color = $hyperlinkColor

    # Latest Handback File (G) - the translated xlf that was handed back
    $gCell = $ws.Cells.Item($Row, 7)
    $gCell.Value = $XlfName
    $ws.Hyperlinks.Add($gCell, $XlfUrl, "", "", $XlfName) | Out-Null
    $gCell.Font.Underline = 2
    $gCell.Font.Color = $hyperlinkColor

    # Latest Handback DateTime (H) - real timestamp instead of 0001-01-01 00:00:00
    $ws.Cells.Item($Row, 8).Value = $HandbackDateTime
}

# ---- zh-cn sheet ----
Set-HandbackRow `
    "zh-cn" `
    2 `
    "2016-03-14 08:32:48" `
    "9c84537b-d831-4af7-88f6-c9357fa0c452.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/c8a43dd5d260f0e80abc4ebdf586354c85331bc9/e2e/9c84537b-d831-4af7-88f6-c9357fa0c452.md" `
    "9c84537b-d831-4af7-88f6-c9357fa0c452.0a5ed60b55d0bb4baed3c8080bc57db20d065c71.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c340b1355a2cd3deba769480323fd5581b40d800/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/9c84537b-d831-4af7-88f6-c9357fa0c452.0a5ed60b55d0bb4baed3c8080bc57db20d065c71.zh-cn.xlf"

Set-HandbackRow `
    "zh-cn" `
    3 `
    "2016-03-14 08:32:48" `
    "d46f061f-73de-4abe-93a9-1c8b5d4dca03.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/c8a43dd5d260f0e80abc4ebdf586354c85331bc9/e2e/d46f061f-73de-4abe-93a9-1c8b5d4dca03.md" `
    "d46f061f-73de-4abe-93a9-1c8b5d4dca03.1b72ec7da7502cf64b3eac4d06b2b442e50818eb.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c340b1355a2cd3deba769480323fd5581b40d800/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/d46f061f-73de-4abe-93a9-1c8b5d4dca03.1b72ec7da7502cf64b3eac4d06b2b442e50818eb.zh-cn.xlf"

# ---- de-de sheet ----
Set-HandbackRow `
    "de-de" `
    2 `
    "2016-03-14 08:32:57" `
    "9c84537b-d831-4af7-88f6-c9357fa0c452.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/c8a43dd5d260f0e80abc4ebdf586354c85331bc9/e2e/9c84537b-d831-4af7-88f6-c9357fa0c452.md" `
    "9c84537b-d831-4af7-88f6-c9357fa0c452.0a5ed60b55d0bb4baed3c8080bc57db20d065c71.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cb644d549d2f80c12f1218faa6d3ef4ed45648c3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/9c84537b-d831-4af7-88f6-c9357fa0c452.0a5ed60b55d0bb4baed3c8080bc57db20d065c71.de-de.xlf"

Set-HandbackRow `
    "de-de" `
    3 `
    "2016-03-14 08:32:57" `
    "d46f061f-73de-4abe-93a9-1c8b5d4dca03.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/c8a43dd5d260f0e80abc4ebdf586354c85331bc9/e2e/d46f061f-73de-4abe-93a9-1c8b5d4dca03.md" `
    "d46f061f-73de-4abe-93a9-1c8b5d4dca03.1b72ec7da7502cf64b3eac4d06b2b442e50818eb.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cb644d549d2f80c12f1218faa6d3ef4ed45648c3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/d46f061f-73de-4abe-93a9-1c8b5d4dca03.1b72ec7da7502cf64b3eac4d06b2b442e50818eb.de-de.xlf"

# ---- Overview sheet ----
# Same "Ready for handoff" -> "Handed back: in sync with en-US" status text,
# mirrored here for both language columns / both files.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 2).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(2, 3).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(3, 2).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"
